$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 99: date/time value changed (A99) ---
$ws.Range("A99").Value = 45478.2916666667

# --- New row 100 ---
# A100: date/time value. Seed the value, then copy A99's formats
# (number format + style) onto it so it gets the same custom date style.
$ws.Range("A100").Value = 45481.6447800926
$ws.Range("A99").Copy()
$ws.Range("A100").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B100").Value = 10800
$ws.Range("C100").Value = 6.26000022888184
$ws.Range("D100").Value = 6.09999990463257
$ws.Range("E100").Value = 6.15999984741211
$ws.Range("F100").Value = 6.1399998664856

# G100 holds the numeric-looking text "6.1399998664856" as a genuine
# string (matches shared string already used elsewhere in the sheet),
# not a number. Force text entry via a temporary text number format,
# then drop the format again so the cell keeps the default style.
$ws.Range("G100").NumberFormat = "@"
$ws.Range("G100").Value = "6.1399998664856"
$ws.Range("G100").ClearFormats()

$ws.Range("H100").Value = "PAL.MI"
